$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old CounterDeal_TC001 row (row 20); this shifts the
# Deals_Chat_* rows (previously 21-23) up to rows 20-22, matching the
# target layout.
$ws.Rows(20).Delete()

# --- Append the CounterDeal_TC001 row back at the end, followed by the
# new CounterDeal / AcceptDeal / DateFilter test rows (rows 23-30).
$newRows = @(
    @(23, "CounterDeal_TC001", "John Tucker", "ONE", "Deal shared successfully"),
    @(24, "CounterDeal_TC002", "Asher Johnson", "ONE", "Deal shared successfully"),
    @(25, "CounterDeal_TC003", "Stan Koster Andersons", "ONE", "Deal shared successfully"),
    @(26, "AcceptDeal_TC001", "John Tucker", "ONE", "Deal shared successfully"),
    @(27, "AcceptDeal_TC002", "Stan Koster Andersons", "ONE", "Deal shared successfully"),
    @(28, "DateFilter_TC001", "John Tucker", "ONE", "Deal shared successfully"),
    @(29, "DateFilter_TC002", "NA", "ALL", "Deal shared successfully"),
    @(30, "DateFilter_TC003", "Stan Koster Andersons", "ONE", "Deal shared successfully")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("A$r").VerticalAlignment = -4108
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}

# --- Match the author's final selection state.
$ws.Range("B21").Select()
